$d = $word.ActiveDocument

$old = "Kampagnendaten 2018 für das Sternbild Perseus: 30. Oktober - 8. November und 29. November - 8. Dezember"
$new = "Kampagnendaten Perseus: 16. bis 25. Januar, 7. bis 16. November, 6. bis 15. Dezember"

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = $old
$rng.Find.Forward = $true
$rng.Find.Wrap = 1

while ($rng.Find.Execute()) {
    # Replace the whole matched range (which spans several differently
    # formatted runs) with a single, plain run: delete the old text and
    # insert the new text as unformatted content.
    $rng.Delete()
    $rng.InsertAfter($new)
    $rng.Collapse(0)
}
